$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# --- Data edit: mark "finished" (column AE) = 1 for every study row (2-101) ---
for ($r = 2; $r -le 101; $r++) {
    $ws.Cells.Item($r, 31).Value = 1
}

# --- View state: scroll/zoom/selection as of the last save ---
$ws.Activate()
$excel.ActiveWindow.ScrollColumn = 11
$excel.ActiveWindow.ScrollRow = 1
$ws.Range("AE105").Select() | Out-Null
$excel.ActiveWindow.Zoom = 134

# --- Workbook window geometry (best effort) ---
$excel.Left = 0
$excel.Top = 740
$excel.Width = 29400
$excel.Height = 16740
